$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(109).Copy()
$ws.Rows.Item(109).Insert()
Write-Host "done"
